$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/enrollment-pcp"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet update ---
# The "Extension" row's Constraint(s) column (AI2) loses its ele-1/ext-1
# constraint text (that text now only applies further down, to the
# Extension.extension row, where it already was unchanged).
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""

# The "Extension.url" row's Fixed Value column (Q5) mirrors the same URL
# string that was updated on the Metadata sheet, so it must be refreshed too.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/enrollment-pcp"
